$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "1")
$ws.Range("B2").Value = 0.8933333333333333
$ws.Range("C2").Value = 0.9241379310344827
$ws.Range("D2").Value = 0.9084745762711864
$ws.Range("E2").Value = 145

# Row 3 (class "2")
$ws.Range("B3").Value = 0.972972972972973
$ws.Range("C3").Value = 0.972972972972973
$ws.Range("D3").Value = 0.972972972972973
$ws.Range("E3").Value = 148

# Row 4 (class "3")
$ws.Range("B4").Value = 0.9363057324840764
$ws.Range("C4").Value = 0.9735099337748344
$ws.Range("D4").Value = 0.9545454545454546
$ws.Range("E4").Value = 151

# Row 5 (class "4")
$ws.Range("B5").Value = 0.9103448275862069
$ws.Range("C5").Value = 0.8461538461538461
$ws.Range("D5").Value = 0.8770764119601329
$ws.Range("E5").Value = 156

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.9283333333333333
$ws.Range("C6").Value = 0.9283333333333333
$ws.Range("D6").Value = 0.9283333333333333
$ws.Range("E6").Value = 0.9283333333333333

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.9282392165941474
$ws.Range("C7").Value = 0.929193670984034
$ws.Range("D7").Value = 0.9282673539374368

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.9282154867364619
$ws.Range("C8").Value = 0.9283333333333333
$ws.Range("D8").Value = 0.9278151624357773

$wb.Save()
